$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20's date (A20) previously used the "last row" date-only format;
# now that it's no longer the last row, give it the regular
# date+time format used by the rest of the series.
$ws.Cells.Item(20, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's results as row 21.
$ws.Cells.Item(21, 1).Value = 45605
$ws.Cells.Item(21, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(21, 2).Value = 48
$ws.Cells.Item(21, 3).Value = 42
$ws.Cells.Item(21, 4).Value = 50
